$d = $word.ActiveDocument
$wmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-RangeXml([int]$paraIndex, [string]$innerRunsXml) {
    $p = $d.Paragraphs($paraIndex)
    $start = $p.Range.Start
    $end = $p.Range.End - 1   # exclude paragraph mark
    $r = $d.Range($start, $end)
    $xml = '<w:document ' + $wmlNs + '><w:body><w:p>' + $innerRunsXml + '</w:p></w:body></w:document>'
    $r.InsertXML($xml)
}

# ------------------------------------------------------------------
# Paragraph 1: Title
# ------------------------------------------------------------------
$d.Content.Find.Execute("Unveiling the Enigma of Time", $true, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "Unraveling the Mysteries of Chemistry: A Journey into the World of Elements and Reactions", 2) | Out-Null

# ------------------------------------------------------------------
# Paragraph 2: Byline  "Madeline Adler" -> "Dr" + "." + " Emily Carter"
# ------------------------------------------------------------------
$byline = '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>Dr</w:t></w:r>' + `
          '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>.</w:t></w:r>' + `
          '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve"> Emily Carter</w:t></w:r>'
Replace-RangeXml 2 $byline

# ------------------------------------------------------------------
# Paragraph 3: Email   "maddie" + "." + "adler@xyzacademy" + "." + "edu"
#                  ->  "emcarter@chemistryeducators" + "." + "org"
# ------------------------------------------------------------------
$email = '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>emcarter@chemistryeducators</w:t></w:r>' + `
         '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>.</w:t></w:r>' + `
         '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>org</w:t></w:r>'
Replace-RangeXml 3 $email

# ------------------------------------------------------------------
# Paragraph 5: Body text
# ------------------------------------------------------------------
$rpr24 = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr>'

$body = ''
$body += '<w:r>' + $rpr24 + '<w:t>Have you ever wondered about the world around you? Why do things change? Why are there so many different substances? The answers to these questions lie in one of the most intriguing and fundamental subjects--chemistry</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:t>.</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:t xml:space="preserve"> On our voyage through the captivating realm of chemistry, we will unravel the secrets of matter, understand how substances interact, and explore the incredible applications of chemistry in our everyday lives</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:t>.</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:br/></w:r>'
$body += '<w:r>' + $rpr24 + '<w:br/><w:t>In this extraordinary odyssey, we will explore the basic building blocks of all matter--the elements</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:t>.</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:t xml:space="preserve"> We will investigate their properties, their bonding behavior, and their arrangements to form different compounds</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:t>.</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:t xml:space="preserve"> Through demonstrations and hands-on experiments, we will uncover the mysteries of chemical reactions, learning how atoms rearrange and energy is transferred</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:t>.</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:t xml:space="preserve"> The world of chemistry is a place of fascinating phenomena, from the colorful fireworks that light up our skies to the complex processes occurring within our bodies</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:t>.</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:br/></w:r>'
$body += '<w:r>' + $rpr24 + '<w:br/><w:t>While embarking on this journey of discovery, we will delve into the diverse applications of chemistry</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:t>.</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:t xml:space="preserve"> From understanding the role of chemistry in fields such as medicine, engineering, and agriculture to comprehending the impact of chemistry on our environment, we will appreciate the significance of this science in shaping our world</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:t>.</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:t xml:space="preserve"> Chemistry provides solutions to real-world problems, enhances our lives, and continues to push the boundaries of human knowledge</w:t></w:r>'
$body += '<w:r>' + $rpr24 + '<w:t>.</w:t></w:r>'

Replace-RangeXml 5 $body

# ------------------------------------------------------------------
# Paragraph 7: Summary
# ------------------------------------------------------------------
$rprS = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr>'

$summary = ''
$summary += '<w:r>' + $rprS + '<w:t>Our exploration of chemistry has unveiled the fundamental concepts and applications of this dynamic science</w:t></w:r>'
$summary += '<w:r>' + $rprS + '<w:t>.</w:t></w:r>'
$summary += '<w:r>' + $rprS + '<w:t xml:space="preserve"> We have journeyed through the world of elements, reactions, and compounds, unraveling the secrets of matter and its transformations</w:t></w:r>'
$summary += '<w:r>' + $rprS + '<w:t>.</w:t></w:r>'
$summary += '<w:r>' + $rprS + '<w:t xml:space="preserve"> Throughout our voyage, we have witnessed the power of chemistry in diverse fields, from medicine to engineering</w:t></w:r>'
$summary += '<w:r>' + $rprS + '<w:t>.</w:t></w:r>'
$summary += '<w:r>' + $rprS + '<w:t xml:space="preserve"> This </w:t></w:r>'
$summary += '<w:r>' + $rprS + '<w:lastRenderedPageBreak/><w:t>knowledge equips us with a deeper understanding of the world around us, empowering us to appreciate the intricacies of chemical processes and their impact on our lives</w:t></w:r>'
$summary += '<w:r>' + $rprS + '<w:t>.</w:t></w:r>'
$summary += '<w:r>' + $rprS + '<w:t xml:space="preserve"> As we continue our scientific odyssey, we look forward to unraveling even more mysteries of the chemical realm</w:t></w:r>'
$summary += '<w:r>' + $rprS + '<w:t>.</w:t></w:r>'

Replace-RangeXml 7 $summary

# ------------------------------------------------------------------
# Add a new empty paragraph at the very end of the body
# ------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$endPt = $p7.Range.End
$insPt = $d.Range($endPt, $endPt)
$insPt.InsertXML('<w:document ' + $wmlNs + '><w:body><w:p/></w:body></w:document>')
